$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Teaching Assistant" -> "Recognised Teacher in Psychology"
# for the three rows it appears in (A10:A12)
$ws.Range("A10").Value = "Recognised Teacher in Psychology"
$ws.Range("A11").Value = "Recognised Teacher in Psychology"
$ws.Range("A12").Value = "Recognised Teacher in Psychology"

# Row 10 and 12 gain the same row height as row 11 (30)
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30

# Update the selection to match the edited range
$ws.Range("A10:A12").Select()
